$wb = $excel.ActiveWorkbook

# Sheet "展览" - update 想去人数 (F column) figures for a few events
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 127
$wsExpo.Range("F4").Value = 172
$wsExpo.Range("F5").Value = 3282
$wsExpo.Range("F7").Value = 13

# Sheet "全部类型" - same events are duplicated here, apply identical updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 127
$wsAll.Range("F4").Value = 172
$wsAll.Range("F5").Value = 3282
$wsAll.Range("F9").Value = 13
